$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add column C containing the same numbers as column B, but stored as text,
# mirroring rows 1-4 (B1:B4 = 1000,2000,3000,4000).
$ws.Range("C1:C4").NumberFormat = "@"
$ws.Range("C1").Value = "1000"
$ws.Range("C2").Value = "2000"
$ws.Range("C3").Value = "3000"
$ws.Range("C4").Value = "4000"

$wb.Save()
